$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.841.12"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.455.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -5.86%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "548.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.26%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.67%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.452.94"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -5.86%  "
$ws.Range("E10").Value = "  -8.00%  "
$ws.Range("E11").Value = "  -1.81%  "
$ws.Range("E12").Value = "  -7.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.351"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.37%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.98"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.898.12"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -5.80%  "
$ws.Range("E16").Value = "  -8.74%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "60.757.47"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -4.35%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.455.05"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -6.19%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.03"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.49%  "
$ws.Range("E20").Value = "  -7.52%  "
$ws.Range("E21").Value = "  -6.99%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "318.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.52%  "
$ws.Range("E23").Value = "  +0.03%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.68%  "
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0₃0976"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -8.50%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.581.61"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.57%  "
$ws.Range("B28").Value = "Bittensor"
$ws.Range("C28").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "542.19"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -6.36%  "
$ws.Range("B29").Value = "Binance-PegBSC-USD"
$ws.Range("C29").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.39%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.57"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.26%  "
$ws.Range("E33").Value = "  -6.55%  "
$ws.Range("E34").Value = "  -7.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.56"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -9.02%  "
$ws.Range("E36").Value = "  -10.72%  "
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "4.83"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.25%  "
$ws.Range("E39").Value = "  -6.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "18.38"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.16%  "
$ws.Range("E41").Value = "  -5.34%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.999"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  -8.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "39.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.20%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.27"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -8.33%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "146.53"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.59%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -7.77%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "20.77"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.96%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0529"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.33%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.582"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -7.45%  "
$ws.Range("E51").Value = "  -6.01%  "
